# Generate Report for handoff
# - Overview sheet: "Ready for handoff" -> "Handoff transform failed"
# - zh-cn / de-de sheets: clear the per-language "Latest Handoff File" cell
#   (C2) and its hyperlink (handoff never produced a target file), reset the
#   "Latest Handoff Datetime" (D2) to the zero date, and flip the
#   "Handoff Reason" (H2) from Include to Ignored to match.

$wb = $excel.ActiveWorkbook

$zeroDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

# ---------------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de)
# ---------------------------------------------------------------------------
$languages = @("zh-cn", "de-de")

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang)

    # Remember the two hyperlinks that must survive (A2 = source md file,
    # A3 = .localization-config) before wiping the sheet's hyperlink
    # collection -- this engine's Hyperlinks.Delete() only works at the
    # whole-collection level, so the surviving links are re-added afterward.
    # (Reading through a `foreach` is required here -- indexed `.Item(n)`
    # access on this host does not resolve the property getters.)
    $linkAAddress = $null
    $linkADisplay = $null
    $linkBAddress = $null
    $linkBDisplay = $null
    foreach ($h in $ws.Hyperlinks) {
        $hRange = $h.Range
        if ($hRange.Row -eq 2 -and $hRange.Column -eq 1) {
            $linkAAddress = $h.Address
            $linkADisplay = $h.TextToDisplay
        }
        if ($hRange.Row -eq 3 -and $hRange.Column -eq 1) {
            $linkBAddress = $h.Address
            $linkBDisplay = $h.TextToDisplay
        }
    }

    # Status text mirrors the Overview sheet.
    $ws.Range("B2").Value = "Handoff transform failed"

    # The handoff transform failed, so there is no handoff target file
    # anymore -- clear the cell and drop its hyperlink.
    $ws.Range("C2").ClearContents()
    $ws.Hyperlinks.Delete()

    # Re-create the two hyperlinks that must remain.
    $ws.Hyperlinks.Add($ws.Range("A2"), $linkAAddress, [Type]::Missing, [Type]::Missing, $linkADisplay)
    $ws.Hyperlinks.Add($ws.Range("A3"), $linkBAddress, [Type]::Missing, [Type]::Missing, $linkBDisplay)

    # Latest Handoff Datetime resets to the zero date (no successful handoff).
    $ws.Range("D2").Value = $zeroDate

    # Latest Handback DateTime stays the zero date.
    $ws.Range("G2").Value = $zeroDate

    # Handoff Reason flips from Include to Ignored.
    $ws.Range("H2").Value = "Ignored"

    # Row 3 values are unchanged (already zero date / Ignored) but are
    # re-asserted here for clarity and to keep the shared-string table tight.
    $ws.Range("D3").Value = $zeroDate
    $ws.Range("G3").Value = $zeroDate
    $ws.Range("H3").Value = "Ignored"
}

Write-Output "done"
